$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "IMAGES" table at J1:K4 -------------------------------------------
# Header J1:K1 = "IMAGES" (bold, centered - same look as the other table headers)
$ws.Range("J1").Value = "IMAGES"
$ws.Range("A1").Copy()
$ws.Range("J1:K1").PasteSpecial(-4122)
$ws.Range("J1:K1").Merge()

# Row 2: id / chiave univoca (same highlight colour as the other "id" rows)
$ws.Range("J2").Value = "id"
$ws.Range("D15").Copy()
$ws.Range("J2").PasteSpecial(-4122)
$ws.Range("K2").Value = "chiave univoca"

# Row 3: name
$ws.Range("J3").Value = "name"

# Row 4: id_user / FK 1 image -> 1 o + user (same highlight colour as the other "id_user" FK row)
$ws.Range("J4").Value = "id_user"
$ws.Range("D20").Copy()
$ws.Range("J4").PasteSpecial(-4122)
$ws.Range("K4").Value = "FK 1 image -> 1 o + user"

$excel.CutCopyMode = $false

# --- Column width tweaks ------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 10.417
$ws.Columns.Item(5).ColumnWidth = 22.085
$ws.Columns.Item(8).ColumnWidth = 68.7541
$ws.Columns.Item(10).ColumnWidth = 6.417
$ws.Columns.Item(11).ColumnWidth = 20.2544

# --- Selection matches the authored workbook -----------------------------
$ws.Range("K6").Select()
